{"js": "// Replace the 25 multiplication-equation answers in the single table of\n// this worksheet. Each entry below is [rowIndex, colIndex, oldText, newText]\n// using 0-based row/column indices into the table (row-major order, same\n// order the cells appear in the document).\nconst replacements = [\n  [0, 0, \"32\u00d728=896\", \"40\u00d748=1920\"],\n  [0, 1, \"62\u00d775=4650\", \"74\u00d798=7252\"],\n  [0, 2, \"25\u00d736=900\", \"94\u00d788=8272\"],\n  [0, 3, \"65\u00d786=5590\", \"85\u00d711=935\"],\n  [0, 4, \"62\u00d779=4898\", \"30\u00d762=1860\"],\n  [4, 0, \"83\u00d731=2573\", \"57\u00d785=4845\"],\n  [4, 1, \"69\u00d727=1863\", \"76\u00d748=3648\"],\n  [4, 2, \"82\u00d774=6068\", \"17\u00d776=1292\"],\n  [4, 3, \"45\u00d799=4455\", \"53\u00d748=2544\"],\n  [4, 4, \"53\u00d759=3127\", \"45\u00d789=4005\"],\n  [9, 0, \"45\u00d785=3825\", \"97\u00d748=4656\"],\n  [9, 1, \"22\u00d798=2156\", \"88\u00d715=1320\"],\n  [9, 2, \"84\u00d729=2436\", \"14\u00d741=574\"],\n  [9, 3, \"54\u00d792=4968\", \"67\u00d722=1474\"],\n  [9, 4, \"96\u00d730=2880\", \"80\u00d793=7440\"],\n  [14, 0, \"26\u00d733=858\", \"42\u00d734=1428\"],\n  [14, 1, \"32\u00d757=1824\", \"62\u00d777=4774\"],\n  [14, 2, \"17\u00d738=646\", \"88\u00d784=7392\"],\n  [14, 3, \"73\u00d760=4380\", \"60\u00d730=1800\"],\n  [14, 4, \"18\u00d783=1494\", \"91\u00d796=8736\"],\n  [19, 0, \"63\u00d796=6048\", \"84\u00d717=1428\"],\n  [19, 1, \"17\u00d715=255\", \"65\u00d714=910\"],\n  [19, 2, \"42\u00d711=462\", \"79\u00d781=6399\"],\n  [19, 3, \"34\u00d788=2992\", \"61\u00d750=3050\"],\n  [19, 4, \"18\u00d783=1494\", \"30\u00d734=1020\"],\n];\n\nconst table = context.document.body.tables.getFirst();\n\nfor (const [rowIndex, colIndex, oldText, newText] of replacements) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const paragraph = cell.body.paragraphs.getFirst();\n  paragraph.load(\"text\");\n  await context.sync();\n\n  if (paragraph.text !== oldText) {\n    throw new Error(\n      `Unexpected cell text at row ${rowIndex}, col ${colIndex}: ` +\n      `expected \"${oldText}\" but found \"${paragraph.text}\"`\n    );\n  }\n\n  // Replace only the text of the run, keeping paragraph/run formatting\n  // (font, size, alignment) untouched.\n  paragraph.getRange().insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 multiplication-equation answers in the single table of\n# this worksheet. Each row below is (tableRow, tableCol, oldText, newText)\n# using 1-based row/column indices, matching Word's Table.Cell(row, col).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @(1, 1, \"32\u00d728=896\", \"40\u00d748=1920\"),\n    @(1, 2, \"62\u00d775=4650\", \"74\u00d798=7252\"),\n    @(1, 3, \"25\u00d736=900\", \"94\u00d788=8272\"),\n    @(1, 4, \"65\u00d786=5590\", \"85\u00d711=935\"),\n    @(1, 5, \"62\u00d779=4898\", \"30\u00d762=1860\"),\n    @(5, 1, \"83\u00d731=2573\", \"57\u00d785=4845\"),\n    @(5, 2, \"69\u00d727=1863\", \"76\u00d748=3648\"),\n    @(5, 3, \"82\u00d774=6068\", \"17\u00d776=1292\"),\n    @(5, 4, \"45\u00d799=4455\", \"53\u00d748=2544\"),\n    @(5, 5, \"53\u00d759=3127\", \"45\u00d789=4005\"),\n    @(10, 1, \"45\u00d785=3825\", \"97\u00d748=4656\"),\n    @(10, 2, \"22\u00d798=2156\", \"88\u00d715=1320\"),\n    @(10, 3, \"84\u00d729=2436\", \"14\u00d741=574\"),\n    @(10, 4, \"54\u00d792=4968\", \"67\u00d722=1474\"),\n    @(10, 5, \"96\u00d730=2880\", \"80\u00d793=7440\"),\n    @(15, 1, \"26\u00d733=858\", \"42\u00d734=1428\"),\n    @(15, 2, \"32\u00d757=1824\", \"62\u00d777=4774\"),\n    @(15, 3, \"17\u00d738=646\", \"88\u00d784=7392\"),\n    @(15, 4, \"73\u00d760=4380\", \"60\u00d730=1800\"),\n    @(15, 5, \"18\u00d783=1494\", \"91\u00d796=8736\"),\n    @(20, 1, \"63\u00d796=6048\", \"84\u00d717=1428\"),\n    @(20, 2, \"17\u00d715=255\", \"65\u00d714=910\"),\n    @(20, 3, \"42\u00d711=462\", \"79\u00d781=6399\"),\n    @(20, 4, \"34\u00d788=2992\", \"61\u00d750=3050\"),\n    @(20, 5, \"18\u00d783=1494\", \"30\u00d734=1020\")\n)\n\nforeach ($entry in $replacements) {\n    $rowIndex = $entry[0]\n    $colIndex = $entry[1]\n    $oldText  = $entry[2]\n    $newText  = $entry[3]\n\n    $cell = $t.Cell($rowIndex, $colIndex)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($current -ne $oldText) {\n        throw \"Unexpected cell text at row ${rowIndex}, col ${colIndex}: expected '$oldText' but found '$current'\"\n    }\n\n    # Assigning .Text keeps the existing run/paragraph formatting\n    # (font, size, alignment) intact - only the characters change.\n    $cell.Range.Text = $newText\n}\n"}
